$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Header row
$ws2.Range("A1").Value = "num"
$ws2.Range("B1").Value = "RC"
$ws2.Range("C1").Value = "XRC"
$ws2.Range("D1").Value = "ERD"

# Data rows (transpose of Sheet1's throughput table)
$ws2.Range("A2").Value = 8
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 6.6
$ws2.Range("D2").Value = 6.6

$ws2.Range("A3").Value = 16
$ws2.Range("B3").Value = 1
$ws2.Range("C3").Value = 7.06
$ws2.Range("D3").Value = 7.06

$ws2.Range("A4").Value = 32
$ws2.Range("B4").Value = 1
$ws2.Range("C4").Value = 7.59
$ws2.Range("D4").Value = 7.59

$ws2.Range("A5").Value = 64
$ws2.Range("B5").Value = 2
$ws2.Range("C5").Value = 8.66
$ws2.Range("D5").Value = 8.66

$ws2.Range("A6").Value = 128
$ws2.Range("B6").Value = 3
$ws2.Range("C6").Value = 11.08
$ws2.Range("D6").Value = 11.08

$ws2.Range("A7").Value = 256
$ws2.Range("B7").Value = 5
$ws2.Range("C7").Value = 14.26
$ws2.Range("D7").Value = 16.63

$ws2.Range("A8").Value = 512
$ws2.Range("B8").Value = 7
$ws2.Range("C8").Value = 24.92
$ws2.Range("D8").Value = 26.83

$ws2.Range("A9").Value = 1024
$ws2.Range("B9").Value = 10
$ws2.Range("C9").Value = 39.67
$ws2.Range("D9").Value = 47.13

$ws2.Range("A10").Value = 2048
$ws2.Range("B10").Value = 17
$ws2.Range("C10").Value = 83.98
$ws2.Range("D10").Value = 93.66

$ws2.Range("A11").Value = 4096
$ws2.Range("B11").Value = 28
$ws2.Range("C11").Value = 92.56
$ws2.Range("D11").Value = 95.71

$ws2.Range("A12").Value = 8192
$ws2.Range("B12").Value = 29
$ws2.Range("C12").Value = 96.74
$ws2.Range("D12").Value = 98.37

$ws2.Range("A13").Value = 16384
$ws2.Range("B13").Value = 33
$ws2.Range("C13").Value = 96.27
$ws2.Range("D13").Value = 97.67

$ws2.Range("A14").Value = 32768
$ws2.Range("B14").Value = 38
$ws2.Range("C14").Value = 95.87
$ws2.Range("D14").Value = 97.07

$ws2.Range("A15").Value = 65536
$ws2.Range("B15").Value = 31
$ws2.Range("C15").Value = 92.9
$ws2.Range("D15").Value = 94.25

$ws2.Range("A16").Value = 131072
$ws2.Range("B16").Value = 29
$ws2.Range("C16").Value = 89.72
$ws2.Range("D16").Value = 90.73

$ws2.Range("A17").Value = 262144
$ws2.Range("B17").Value = 37
$ws2.Range("C17").Value = 89.32
$ws2.Range("D17").Value = 89.87

$ws2.Range("A18").Value = 524288
$ws2.Range("B18").Value = 30
$ws2.Range("C18").Value = 78.48
$ws2.Range("D18").Value = 79.06

$ws2.Range("A19").Value = 1048576
$ws2.Range("B19").Value = 30
$ws2.Range("C19").Value = 78.34
$ws2.Range("D19").Value = 79.06

# Select the source data range on Sheet1 (as if it had just been copied from there)
$ws1.Activate() | Out-Null
$ws1.Range("A3:S6").Select() | Out-Null

# Make Sheet2 the active tab to match the saved view state (tabSelected / activeTab)
$ws2.Activate() | Out-Null
